# Final daily scrum update, other doc work
#
# This script brings the "Team Member Report" workbook up to date for the
# daily scrum: it fills in each team member's "next sprint role" cell
# (previously a "TDB" placeholder), records Arpit's completed work for the
# current sprint, and re-points the active sheet/selection back to Arpit.

$wb = $excel.ActiveWorkbook

$arpit    = $wb.Worksheets.Item("Arpit")
$brody    = $wb.Worksheets.Item("Brody")
$michael  = $wb.Worksheets.Item("Michael")
$sakshyam = $wb.Worksheets.Item("Sakshyam")
$vasilis  = $wb.Worksheets.Item("Vasilis")
$yong     = $wb.Worksheets.Item("Yong")

# --- Vasilis: role for next sprint (new shared string, write first so it
#     lands before the other brand-new strings introduced below) ---
$vasilis.Range("B8").Value = " Product Owner, Development"

# --- Arpit: sprint wrap-up (work done this week / issues) ---
$arpit.Range("B7").Value  = "No work completed"
$arpit.Range("B10").Value = "No issues resolved"
$arpit.Range("B9").Value  = "N/A"

# --- Arpit: role for next sprint ---
$arpit.Range("B8").Value = "Development Team"

# --- Remaining team members: role for next sprint ---
$brody.Range("B8").Value    = "Development Team"
$michael.Range("B8").Value  = "Development Team"
$sakshyam.Range("B8").Value = "Scrum Master, Development Team"
$yong.Range("B8").Value     = "Development Team"

# --- View state: Michael's sheet scrolled/zoomed/selected differently ---
$michael.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.Zoom = 189
$michael.Range("C8").Select()

# --- View state: Yong's sheet zoom stays, selection moves, no longer the
#     active/selected tab ---
$yong.Activate()
$excel.ActiveWindow.Zoom = 175
$yong.Range("C8").Select()

# --- Finally, Arpit becomes the active sheet again with its own zoom and
#     selection (this also makes Arpit's tab the one marked selected, and
#     clears it from Yong) ---
$arpit.Activate()
$excel.ActiveWindow.Zoom = 173
$arpit.Range("B15").Select()
